# Enigma workbook — the wiring table on Sheet2 is being worked on next, so
# the author flips the active tab from Sheet3 over to Sheet2 and leaves the
# selection just past the wiring table (I27). Sheet3 keeps its old
# selection (M26) but is no longer the tab that's active/highlighted.
# Recalculating also refreshes the volatile RAND() driver column (A1:A26 on
# Sheet2) with new cached values, as happens on every open/recalc.

$wb = $excel.ActiveWorkbook

$sheet2 = $wb.Worksheets.Item("Sheet2")
$sheet3 = $wb.Worksheets.Item("Sheet3")

# Recalculate everything (including the volatile RAND() column) so the
# cached <v> values get refreshed, same as Excel does on every recalc.
$null = $excel.CalculateFull()

# Sheet3 stays parked on M26, it's just not the active tab anymore.
$null = $sheet3.Activate()
$null = $sheet3.Range("M26").Select()

# Sheet2 becomes the active tab, selection lands just below the wiring
# table instead of the old H1:I26 block.
$null = $sheet2.Activate()
$null = $sheet2.Range("I27").Select()
